$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# New data rows (row => Hinge Direction (A), Door Type (C), Width Decimal (D), Height Decimal (F), Quantity (H))
$data = @{
    2  = @("Left",  "Single",     33, 77, 2)
    3  = @("Right", "Lazy Susan", 22, 88, 3)
    4  = @("Left",  "Pair",       23, 55, 4)
    5  = @("Right", "Lazy Susan", 34, 32, 5)
    6  = @("Left",  "Single",     12, 55, 6)
    7  = @("Right", "Lazy Susan", 40, 55, 7)
    8  = @("Left",  "Single",     12, 77, 8)
    9  = @("Right", "Lazy Susan", 41, 42, 4)
    10 = @("Left",  "Pair",       33, 44, 1)
    11 = @("Right", "Single",      5, 99, 2)
    12 = @("Left",  "Lazy Susan",  9,  5, 4)
    13 = @("Right", "Single",     33,  6, 1)
    14 = @("Left",  "Lazy Susan", 27, 12, 6)
    15 = @("Right", "Single",     11,  5, 7)
    16 = @("Left",  "Lazy Susan", 14, 11, 8)
    17 = @("Right", "Lazy Susan", 16, 26, 9)
    18 = @("Left",  "Single",     18, 28, 2)
    19 = @("Right", "Lazy Susan", 29, 30, 4)
    20 = @("Left",  "Single",      6, 37, 4)
}

$rows = 2..20

# Update column by column (C, then D, then F, then H, then A) so that new
# shared-string entries are appended in the same order the source workbook used:
# Single, Pair, Left.
foreach ($r in $rows) { $ws.Range("C$r").Value = $data[$r][1] }
foreach ($r in $rows) { $ws.Range("D$r").Value = $data[$r][2] }
foreach ($r in $rows) { $ws.Range("F$r").Value = $data[$r][3] }
foreach ($r in $rows) { $ws.Range("H$r").Value = $data[$r][4] }
foreach ($r in $rows) { $ws.Range("A$r").Value = $data[$r][0] }

# Update sheet view (top-left cell + selection)
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F21").Select()
